# Eduati SW620_noCTRL_meas.xlsx - "bug fix in Eduati data files"
#
# Sheet1 ("Condition" index sheet) had stray leftover index values in
# column A for rows 45:87 (no corresponding data in columns B:N - the
# real data only goes through row 44, matching Sheet2/Sheet3). Remove
# those extra rows, and make Sheet1 the active/selected sheet (it had
# been left on Sheet3).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the 43 leftover rows (45:87) that only held stray index numbers
# in column A with no real measurement data.
$ws1.Rows("45:87").Delete() | Out-Null

# Sheet1 should be the active tab/selection when the workbook is
# reopened (previously Sheet3 was left tabSelected/active).
$ws1.Activate() | Out-Null
$ws1.Range("G58").Select() | Out-Null
